$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the shared formula across E7:E8 first (before the row delete below)
# so the engine keeps them grouped as a single shared formula, matching the
# original authoring (D7/C7 shared across the range).
$ws.Range("E7:E8").Formula = "=D7/C7"

# Update iteration 6 (row 8) task counts: 12 -> 17 for both estimated and actual
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 17

# Update the "Action Done" text for iteration 5 (row 7) and iteration 6 (row 8)
# to reflect that there were no unplanned tasks.
$newText = "Our estimates are fairly accurate, the team are roughly on track. No unplanned tasks. No spillovers."
$ws.Range("F7").Value = $newText
$ws.Range("F8").Value = $newText

# Remove the now-empty/erroring iteration 7 row (row 9) entirely.
$ws.Rows.Item(9).Delete()

# Update selection to match the post-edit state (cursor moved to row 9, full row selected)
$ws.Range("A9:XFD9").Select()
